$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update category label in B2 from "Autoradios" to "Departamentos"
$ws.Range("B2").Value = "Departamentos"

# Clear the contents of row 3 (A3:C3) while keeping formatting
$ws.Range("A3:C3").ClearContents()

# Update the active selection to G8
$ws.Range("G8").Select()
